$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text for column A (unit changed from bar to MPa)
$ws.Range("A1").Value = "LH2 Feed pressure (MPa)"

# New simulation results for 0.5s timestep run
$data = @(
    @(0.01,                   32.592550571366232, 15.990458615771146, 0.49061697643937568, 30.794038459236532, 15.234835784775994, 0.49473328433173969, 385.5, 353),
    @(0.020000000000000004,   32.932372017592201, 16.078444421406854, 0.48822612634212564, 30.117062545299198, 15.296917341733348, 0.50791531606793294, 322.5, 288),
    @(0.030000000000000006,   33.1796975239291,   16.217333375594809, 0.48877279137035279, 30.194228218778868, 15.288582615456379, 0.50634122868382725, 297,   261),
    @(0.040000000000000001,   33.39462924544231,  16.317367183436009, 0.48862249865112661, 30.590676489596902, 15.480958121910632, 0.50606785787085506, 281.5, 247),
    @(0.05000000000000001,    33.588616560882791, 16.542366296738358, 0.49249918545331073, 30.634043816811094, 15.356483311938568, 0.50128815522263392, 271,   237.5),
    @(0.060000000000000012,   34.038245924791944, 16.776145065184295, 0.49286162107916665, 30.907467922690074, 15.473842353435121, 0.50065059978838711, 264,   231.5),
    @(0.070000000000000007,   34.506021199950069, 16.966938390040298, 0.49170949880668507, 30.940195223590329, 15.498883332613644, 0.50093036648962486, 258,   225.5),
    @(0.080000000000000002,   34.542526943152097, 16.836710504553999, 0.48741976903608675, 30.858930562263019, 15.463301739826747, 0.50109648837722898, 253,   221),
    @(0.090000000000000011,   35.218313745201684, 17.311994536218648, 0.49156227812233966, 31.060537986029125, 15.553497301966893, 0.50074783987845861, 249.5, 217.5),
    @(0.10000000000000001,    35.045649540805663, 17.234527415253048, 0.49177366209708562, 31.13447520467524,  15.599402815725494, 0.5010331060079356,  245.5, 214.5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
}

# Auto-fit columns to reflect the new (slightly wider) content widths
$newColumnWidths = @(23, 33.28515625, 31.85546875, 38.7109375, 32, 30.5703125, 37.42578125, 32.28515625, 31)
for ($c = 0; $c -lt $newColumnWidths.Length; $c++) {
    $ws.Columns.Item($c + 1).ColumnWidth = $newColumnWidths[$c] - 0.8333333333333334
}
